$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New attendee checked in: replace the old sample row with the latest scan.
# Student ID keeps leading zeros, so force text formatting before typing it.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "0003694140"

# Name wasn't captured for this scan -> leave it blank (still a text cell).
$ws.Range("B2").Formula = "=""""" 

$ws.Range("C2").Value = "15:06:14"

# Move the selection back to A1 to match the saved view state.
$ws.Range("A1").Select()
